# [ADD] Schedule tasks for today's day of the week only
#
# Rearranges the "meetings" sheet columns to (day, time, meeting_id,
# passcode, type, user), consolidates both meetings onto the same day
# ("friday") and same time, renames the old "title" header to "user",
# and drops the now-unused "monday" row / trailing blank rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# ---------------------------------------------------------------
# 1) Stage the existing cell *styles* we still need into scratch
#    cells far away from the used range, so that re-using a style
#    later isn't affected by us having already overwritten its
#    original source cell.
# ---------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("Z1").PasteSpecial($xlPasteFormats)   # s1 - bold header, centered+wrap
$ws.Range("E1").Copy()
$ws.Range("Z2").PasteSpecial($xlPasteFormats)   # s2 - bold header
$ws.Range("A2").Copy()
$ws.Range("Z3").PasteSpecial($xlPasteFormats)   # s3 - body, vcenter+wrap
$ws.Range("D2").Copy()
$ws.Range("Z4").PasteSpecial($xlPasteFormats)   # s4 - time, vcenter+wrap
$ws.Range("B10").Copy()
$ws.Range("Z5").PasteSpecial($xlPasteFormats)   # s5 - underline font
$ws.Range("A3").Copy()
$ws.Range("Z6").PasteSpecial($xlPasteFormats)   # s6 - hyperlink style
$excel.CutCopyMode = $false

function Set-FormatFrom($srcAddr, $destAddr) {
    $ws.Range($srcAddr).Copy()
    $ws.Range($destAddr).PasteSpecial($xlPasteFormats)
}

# ---------------------------------------------------------------
# 2) Clear everything in the data area first so stale values/styles
#    from the old layout don't linger in cells we don't explicitly
#    touch afterwards.
# ---------------------------------------------------------------
$ws.Range("A1:F10").Clear()

# ---------------------------------------------------------------
# 3) Header row (row 1) - style s1 for A:D, s2 for E:F
# ---------------------------------------------------------------
foreach ($col in @("A1","B1","C1","D1")) { Set-FormatFrom "Z1" $col }
foreach ($col in @("E1","F1"))           { Set-FormatFrom "Z2" $col }

$ws.Range("A1").Value = "day"
$ws.Range("B1").Value = "time"
$ws.Range("C1").Value = "meeting_id"
$ws.Range("D1").Value = "passcode"
$ws.Range("E1").Value = "type"
$ws.Range("F1").Value = "user"

# ---------------------------------------------------------------
# 4) Row 2 - juan's google_meet entry
# ---------------------------------------------------------------
Set-FormatFrom "Z3" "A2"
Set-FormatFrom "Z4" "B2"
Set-FormatFrom "Z6" "C2"
Set-FormatFrom "Z3" "D2"

$ws.Range("A2").Value = "friday"
$ws.Range("B2").Value = 0.64722222222222225
$ws.Range("C2").Value = "https://meet.google.com/url_to_your_metting"
$ws.Range("D2").ClearContents()
$ws.Range("E2").Value = "google_meet"
$ws.Range("F2").Value = "juan"

# ---------------------------------------------------------------
# 5) Row 3 - pedro's zoom entry
# ---------------------------------------------------------------
Set-FormatFrom "Z3" "A3"
Set-FormatFrom "Z4" "B3"
Set-FormatFrom "Z3" "C3"
Set-FormatFrom "Z3" "D3"

$ws.Range("A3").Value = "friday"
$ws.Range("B3").Value = 0.64722222222222225
$ws.Range("C3").Value = 1234567890
$ws.Range("D3").Value = 123456
$ws.Range("E3").Value = "zoom"
$ws.Range("F3").Value = "pedro"

# ---------------------------------------------------------------
# 6) Rows 4-7 - empty (formatted) rows
# ---------------------------------------------------------------
foreach ($r in 4..7) {
    Set-FormatFrom "Z3" "A$r"
    Set-FormatFrom "Z4" "B$r"
}
Set-FormatFrom "Z6" "C4"
Set-FormatFrom "Z3" "C5"
Set-FormatFrom "Z6" "C6"
Set-FormatFrom "Z3" "C7"
Set-FormatFrom "Z3" "D4"
Set-FormatFrom "Z3" "D5"

# ---------------------------------------------------------------
# 7) Row 9 - trailing underline-styled row
# ---------------------------------------------------------------
Set-FormatFrom "Z5" "A9"
Set-FormatFrom "Z5" "D9"

# ---------------------------------------------------------------
# 8) Clean up scratch cells used for style staging
# ---------------------------------------------------------------
$ws.Range("Z1:Z6").Clear()
$excel.CutCopyMode = $false

# ---------------------------------------------------------------
# 9) Data validation: the "day" list now lives in column A
# ---------------------------------------------------------------
$ws.Range("A1:A1048576").Validation.Delete()
$ws.Range("A1:A1048576").Validation.Add(3, 1, 1, '"monday,tuesday,wednesday,thursday,friday,saturday,sunday"')
$ws.Range("A1:A1048576").Validation.InputMessage = ""
$ws.Range("A1:A1048576").Validation.ErrorMessage = ""

# ---------------------------------------------------------------
# 10) Hyperlink now anchors on C2 (the meeting_id column)
# ---------------------------------------------------------------
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C2"), "https://meet.google.com/url_to_your_metting") | Out-Null

# ---------------------------------------------------------------
# 11) Column widths follow the new column order
# ---------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 24.33203125
$ws.Columns.Item(2).ColumnWidth = 23.6640625
$ws.Columns.Item(3).ColumnWidth = 41
$ws.Columns.Item(4).ColumnWidth = 24.33203125

# ---------------------------------------------------------------
# 12) Active selection as left by the edit
# ---------------------------------------------------------------
$ws.Range("B7").Select()

Write-Output "applied"
